# Update "想去人数" (number of people interested) values for a few events
# across the "展览" and "全部类型" worksheets, per the new data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3797
$ws1.Range("F6").Value = 163
$ws1.Range("F8").Value = 218
$ws1.Range("F9").Value = 5

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 3797
$ws4.Range("F10").Value = 163
$ws4.Range("F13").Value = 218
$ws4.Range("F14").Value = 5
